# Applies the updated crypto price/volume data (and the USDe/Bittensor row-order swap)
# to Sheet1, preserving each cell as plain text (matches the original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. "1.00") are not
    # auto-converted to numbers, then drop the temporary format so the cell keeps
    # its original (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 'D2' '61.693.52'
Set-TextValue 'E2' '  -1.91%  '
Set-TextValue 'D3' '2.898.76'
Set-TextValue 'E3' '  -1.76%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '567.32'
Set-TextValue 'E5' '  -4.21%  '
Set-TextValue 'D6' '143.38'
Set-TextValue 'E6' '  -3.13%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.502'
Set-TextValue 'E8' '  -0.48%  '
Set-TextValue 'D9' '2.895.34'
Set-TextValue 'E9' '  -1.93%  '
Set-TextValue 'D10' '6.92'
Set-TextValue 'E10' '  -2.46%  '
Set-TextValue 'E11' '  -2.14%  '
Set-TextValue 'D12' '0.430'
Set-TextValue 'E12' '  -2.19%  '
Set-TextValue 'D13' '0.0000231'
Set-TextValue 'E13' '  -1.00%  '
Set-TextValue 'D14' '31.80'
Set-TextValue 'E14' '  -2.64%  '
Set-TextValue 'E15' '  -0.45%  '
Set-TextValue 'D16' '3.376.43'
Set-TextValue 'E16' '  -1.91%  '
Set-TextValue 'D17' '61.648.26'
Set-TextValue 'E17' '  -2.04%  '
Set-TextValue 'D18' '6.55'
Set-TextValue 'E18' '  -1.79%  '
Set-TextValue 'D19' '2.896.24'
Set-TextValue 'E19' '  -1.94%  '
Set-TextValue 'D20' '431.45'
Set-TextValue 'E20' '  -2.05%  '
Set-TextValue 'D21' '13.05'
Set-TextValue 'E21' '  -3.20%  '
Set-TextValue 'D22' '0.654'
Set-TextValue 'E22' '  -2.06%  '
Set-TextValue 'D23' '6.83'
Set-TextValue 'E23' '  -2.55%  '
Set-TextValue 'D24' '79.16'
Set-TextValue 'E24' '  -2.04%  '
Set-TextValue 'D25' '11.89'
Set-TextValue 'E25' '  +0.69%  '
Set-TextValue 'E26' '  +0.04%  '
Set-TextValue 'D27' '9.95'
Set-TextValue 'E27' '  -11.26%  '
Set-TextValue 'E28' '  -5.51%  '
Set-TextValue 'E29' '  +3.67%  '
Set-TextValue 'D30' '7.01'
Set-TextValue 'E30' '  -4.05%  '
Set-TextValue 'E31' '  -4.35%  '
Set-TextValue 'E32' '  -8.61%  '
Set-TextValue 'D33' '0.998'
Set-TextValue 'E33' '  -0.22%  '
Set-TextValue 'E34' '  -1.92%  '
Set-TextValue 'D35' '25.52'
Set-TextValue 'E35' '  -3.19%  '
Set-TextValue 'D36' '0.957'
Set-TextValue 'E36' '  -3.51%  '
Set-TextValue 'D37' '5.37'
Set-TextValue 'E37' '  -4.38%  '
Set-TextValue 'D38' '48.85'
Set-TextValue 'E38' '  -1.76%  '
Set-TextValue 'E39' '  -5.08%  '
Set-TextValue 'D40' '2.81'
Set-TextValue 'E40' '  -8.01%  '
Set-TextValue 'E41' '  -3.44%  '
Set-TextValue 'E42' '  -3.64%  '
Set-TextValue 'D43' '39.81'
Set-TextValue 'E43' '  +0.60%  '
Set-TextValue 'E44' '  -4.05%  '
Set-TextValue 'D45' '2.686.76'
Set-TextValue 'E45' '  -0.82%  '
Set-TextValue 'D46' '132.28'
Set-TextValue 'E46' '  -2.54%  '
Set-TextValue 'D47' '0.0335'
Set-TextValue 'E47' '  -1.06%  '
Set-TextValue 'B48' 'Bittensor'
Set-TextValue 'C48' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D48' '342.86'
Set-TextValue 'E48' '  -4.70%  '
Set-TextValue 'B49' 'USDe'
Set-TextValue 'C49' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D49' '1.00'
Set-TextValue 'E49' '  -0.02%  '
Set-TextValue 'D50' '0.103'
Set-TextValue 'E50' '  -1.66%  '
Set-TextValue 'D51' '21.52'
Set-TextValue 'E51' '  -5.30%  '
